$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table title: roll the reporting window forward one month (October -> November 2016)
$ws.Range("A1").Value = "Table 6.7.A. Capacity Factors for Utility Scale Generators Primarily Using Fossil Fuels, January 2013-November 2016"

# Make room for a new "November" data row just above the footnote row (currently row 45),
# shifting the footnote row (and its formatting) down to row 46.
$ws.Range("A45:I45").Insert(-4121)

# Clone the formatting of the previous month's row (October, row 44) into the new row 45
$ws.Range("A44:I44").Copy()
$ws.Range("A45:I45").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill in the November monthly capacity factors
$ws.Range("A45").Value = "November"
$ws.Range("B45").Value = 0.455
$ws.Range("C45").Value = 0.469
$ws.Range("D45").Value = 0.066
$ws.Range("E45").Value = 0.059
$ws.Range("F45").Value = "NA"
$ws.Range("G45").Value = 0.093
$ws.Range("H45").Value = 0.007
$ws.Range("I45").Value = "NA"
